$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Requirement_Traceability_Matrix")
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(13).Copy()
$ws.Rows.Item(14).PasteSpecial(-4122)
Write-Host "done"
